$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "N native speakers of Spanish participated..." -> "Forty nine native
#    speakers of Spanish participated..."
#    The "N" sits at the very start of its own run/paragraph, so deleting it
#    and inserting the replacement text right before the remaining text keeps
#    the two pieces as separate runs (matching the authored edit, which split
#    the original run into two runs).
# ---------------------------------------------------------------------------
$findR = $d.Content
$findR.Find.Execute("N native speakers of Spanish participated", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$nStart = $findR.Start

$nRange = $d.Range($nStart, $nStart + 1)
$nRange.Delete()
$insertPoint = $d.Range($nStart, $nStart)
$insertPoint.InsertBefore("Forty nine")

# ---------------------------------------------------------------------------
# 2) "age range = 18-35yo)" -> "age range = 20-35yo)"
#    "18" sits in the middle of a run, so a plain delete+insert at that point
#    gets silently re-merged with its (identically formatted) neighbours.
#    Toggling a character property around the edit forces the engine to keep
#    the edited text as its own run, reproducing the three-way run split
#    ("= ", "20", "-35yo...") seen in the authored edit.
# ---------------------------------------------------------------------------
$findR2 = $d.Content
$findR2.Find.Execute("= 18-35yo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$eqStart = $findR2.Start

$yearRange = $d.Range($eqStart + 2, $eqStart + 4)
$yearRange.Font.Bold = $true
$yearRange.Text = "20"
$yearRange.Font.Bold = $false

# ---------------------------------------------------------------------------
# 3) Fix typo in comment #12 (Katerina Tsaroucha's note about trial duration):
#    "...but it seems weird to have there two different numbers there." ->
#    "...but it seems weird to have two different numbers there."
# ---------------------------------------------------------------------------
$comment = $d.Comments.Item(12)
$comment.Range.Text = "The total duration of the trial was different that the time the participants were given to respond, because of the technical issue we had to solve with the microphone. I am not sure it should be mentioned, but it seems weird to have two different numbers there."
